$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 2
$ws.Range("D2").Value = 1

$ws.Range("B5").Value = 0.6666666666666666
$ws.Range("D5").Value = 0.3333333333333333
